# Applies the scheduled-runner price/profit refresh to the Leve tables.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC!row15
$wsALC.Cells.Item(15, 8).Value = 64.73
$wsALC.Cells.Item(15, 9).Value = 64.73
$wsALC.Cells.Item(15, 11).Value = 194.19
$wsALC.Cells.Item(15, 13).Value = -25.19

# ALC!row17
$wsALC.Cells.Item(17, 8).Value = 440759
$wsALC.Cells.Item(17, 10).Value = 440759
$wsALC.Cells.Item(17, 12).Value = 1322277
$wsALC.Cells.Item(17, 14).Value = -1322613

# ALC!row33
$wsALC.Cells.Item(33, 8).Value = 313.375
$wsALC.Cells.Item(33, 9).Value = 172.75
$wsALC.Cells.Item(33, 10).Value = 360.25
$wsALC.Cells.Item(33, 11).Value = 172.75
$wsALC.Cells.Item(33, 12).Value = 360.25
$wsALC.Cells.Item(33, 13).Value = 56.25
$wsALC.Cells.Item(33, 14).Value = -818.25

# ALC!row76
$wsALC.Cells.Item(76, 8).Value = 65220130
$wsALC.Cells.Item(76, 9).Value = 71431280
$wsALC.Cells.Item(76, 10).Value = 2999.5
$wsALC.Cells.Item(76, 11).Value = 71431280
$wsALC.Cells.Item(76, 12).Value = 2999.5
$wsALC.Cells.Item(76, 13).Value = -71430965
$wsALC.Cells.Item(76, 14).Value = -3629.5

# ALC!row79
$wsALC.Cells.Item(79, 8).Value = 65220130
$wsALC.Cells.Item(79, 9).Value = 71431280
$wsALC.Cells.Item(79, 10).Value = 2999.5
$wsALC.Cells.Item(79, 11).Value = 71431280
$wsALC.Cells.Item(79, 12).Value = 2999.5
$wsALC.Cells.Item(79, 13).Value = -71430188
$wsALC.Cells.Item(79, 14).Value = -5183.5

# ALC!row138
$wsALC.Cells.Item(138, 8).Value = 2780.42
$wsALC.Cells.Item(138, 9).Value = 1152.5128
$wsALC.Cells.Item(138, 10).Value = 3821.2131
$wsALC.Cells.Item(138, 11).Value = 3457.5384
$wsALC.Cells.Item(138, 12).Value = 11463.6393
$wsALC.Cells.Item(138, 13).Value = 1682.4616
$wsALC.Cells.Item(138, 14).Value = -21743.6393

# ARM!row19
$wsARM.Cells.Item(19, 8).Value = 2000
$wsARM.Cells.Item(19, 10).Value = 2000
$wsARM.Cells.Item(19, 12).Value = 2000
$wsARM.Cells.Item(19, 14).Value = -2458

# ARM!row32
$wsARM.Cells.Item(32, 8).Value = 16148.029
$wsARM.Cells.Item(32, 9).Value = 10390.113
$wsARM.Cells.Item(32, 10).Value = 67146.71000000001
$wsARM.Cells.Item(32, 11).Value = 10390.113
$wsARM.Cells.Item(32, 12).Value = 67146.71000000001
$wsARM.Cells.Item(32, 13).Value = -10103.113
$wsARM.Cells.Item(32, 14).Value = -67720.71000000001

# ARM!row63
$wsARM.Cells.Item(63, 8).Value = 2097.1875
$wsARM.Cells.Item(63, 9).Value = 2097.1875
$wsARM.Cells.Item(63, 10).Value = 0
$wsARM.Cells.Item(63, 11).Value = 2097.1875
$wsARM.Cells.Item(63, 12).Value = 0
$wsARM.Cells.Item(63, 13).Value = -1411.1875
$wsARM.Cells.Item(63, 14).ClearContents()

# ARM!row66
$wsARM.Cells.Item(66, 8).Value = 2097.1875
$wsARM.Cells.Item(66, 9).Value = 2097.1875
$wsARM.Cells.Item(66, 10).Value = 0
$wsARM.Cells.Item(66, 11).Value = 10485.9375
$wsARM.Cells.Item(66, 12).Value = 0
$wsARM.Cells.Item(66, 13).Value = -7053.9375
$wsARM.Cells.Item(66, 14).ClearContents()

# ARM!row74
$wsARM.Cells.Item(74, 8).Value = 830.7213
$wsARM.Cells.Item(74, 9).Value = 771.25
$wsARM.Cells.Item(74, 10).Value = 1050.3077
$wsARM.Cells.Item(74, 11).Value = 771.25
$wsARM.Cells.Item(74, 12).Value = 1050.3077
$wsARM.Cells.Item(74, 13).Value = 102.75
$wsARM.Cells.Item(74, 14).Value = -2798.3077

# ARM!row77
$wsARM.Cells.Item(77, 8).Value = 830.7213
$wsARM.Cells.Item(77, 9).Value = 771.25
$wsARM.Cells.Item(77, 10).Value = 1050.3077
$wsARM.Cells.Item(77, 11).Value = 3856.25
$wsARM.Cells.Item(77, 12).Value = 5251.538500000001
$wsARM.Cells.Item(77, 13).Value = 511.75
$wsARM.Cells.Item(77, 14).Value = -13987.5385

# BSM!row99
$wsBSM.Cells.Item(99, 8).Value = 17859214
$wsBSM.Cells.Item(99, 9).Value = 38463616
$wsBSM.Cells.Item(99, 10).Value = 2066.6667
$wsBSM.Cells.Item(99, 11).Value = 38463616
$wsBSM.Cells.Item(99, 12).Value = 2066.6667
$wsBSM.Cells.Item(99, 13).Value = -38462118
$wsBSM.Cells.Item(99, 14).Value = -5062.6667

# BSM!row105
$wsBSM.Cells.Item(105, 8).Value = 3376.8635
$wsBSM.Cells.Item(105, 9).Value = 4247.4165
$wsBSM.Cells.Item(105, 10).Value = 2332.2
$wsBSM.Cells.Item(105, 11).Value = 4247.4165
$wsBSM.Cells.Item(105, 12).Value = 2332.2
$wsBSM.Cells.Item(105, 13).Value = -2500.4165
$wsBSM.Cells.Item(105, 14).Value = -5826.2

# CRP!row134
$wsCRP.Cells.Item(134, 8).Value = 20001224
$wsCRP.Cells.Item(134, 9).Value = 1300.6818
$wsCRP.Cells.Item(134, 10).Value = 166667330
$wsCRP.Cells.Item(134, 11).Value = 3902.0454
$wsCRP.Cells.Item(134, 12).Value = 500001990
$wsCRP.Cells.Item(134, 13).Value = -1367.0454
$wsCRP.Cells.Item(134, 14).Value = -500007060

# CUL!row5
$wsCUL.Cells.Item(5, 8).Value = 2226.2292
$wsCUL.Cells.Item(5, 10).Value = 2495
$wsCUL.Cells.Item(5, 12).Value = 7485
$wsCUL.Cells.Item(5, 14).Value = -7709

# CUL!row50
$wsCUL.Cells.Item(50, 8).Value = 222.5
$wsCUL.Cells.Item(50, 9).Value = 230
$wsCUL.Cells.Item(50, 11).Value = 690
$wsCUL.Cells.Item(50, 13).Value = -209

# CUL!row53
$wsCUL.Cells.Item(53, 8).Value = 222.5
$wsCUL.Cells.Item(53, 9).Value = 230
$wsCUL.Cells.Item(53, 11).Value = 690
$wsCUL.Cells.Item(53, 13).Value = -209

# CUL!row114
$wsCUL.Cells.Item(114, 8).Value = 1754.2142
$wsCUL.Cells.Item(114, 9).Value = 500
$wsCUL.Cells.Item(114, 10).Value = 1850.6923
$wsCUL.Cells.Item(114, 11).Value = 1500
$wsCUL.Cells.Item(114, 12).Value = 5552.0769
$wsCUL.Cells.Item(114, 13).Value = 1754
$wsCUL.Cells.Item(114, 14).Value = -12060.0769

# CUL!row132
$wsCUL.Cells.Item(132, 8).Value = 1010.0455
$wsCUL.Cells.Item(132, 9).Value = 932.875
$wsCUL.Cells.Item(132, 10).Value = 1215.8334
$wsCUL.Cells.Item(132, 11).Value = 8395.875
$wsCUL.Cells.Item(132, 12).Value = 10942.5006
$wsCUL.Cells.Item(132, 13).Value = -5865.875
$wsCUL.Cells.Item(132, 14).Value = -16002.5006

# CUL!row135
$wsCUL.Cells.Item(135, 8).Value = 2226.2292
$wsCUL.Cells.Item(135, 10).Value = 2495
$wsCUL.Cells.Item(135, 12).Value = 22455
$wsCUL.Cells.Item(135, 14).Value = -27525

# CUL!row137
$wsCUL.Cells.Item(137, 8).Value = 51172.047
$wsCUL.Cells.Item(137, 9).Value = 2931.7693
$wsCUL.Cells.Item(137, 10).Value = 129562.5
$wsCUL.Cells.Item(137, 11).Value = 8795.3079
$wsCUL.Cells.Item(137, 12).Value = 388687.5
$wsCUL.Cells.Item(137, 13).Value = -3695.3079
$wsCUL.Cells.Item(137, 14).Value = -398887.5

# GSM!row18
$wsGSM.Cells.Item(18, 8).Value = 2053162.4
$wsGSM.Cells.Item(18, 10).Value = 66453
$wsGSM.Cells.Item(18, 12).Value = 66453
$wsGSM.Cells.Item(18, 14).Value = -67039

# GSM!row126
$wsGSM.Cells.Item(126, 8).Value = 3080.4443
$wsGSM.Cells.Item(126, 9).Value = 3080.4443
$wsGSM.Cells.Item(126, 11).Value = 9241.332900000001
$wsGSM.Cells.Item(126, 13).Value = -6771.332900000001

# GSM!row132
$wsGSM.Cells.Item(132, 8).Value = 2883.5417
$wsGSM.Cells.Item(132, 9).Value = 2660.3057
$wsGSM.Cells.Item(132, 11).Value = 7980.9171
$wsGSM.Cells.Item(132, 13).Value = -5450.9171

# LTW!row68
$wsLTW.Cells.Item(68, 8).Value = 2639.3462
$wsLTW.Cells.Item(68, 10).Value = 2862.5
$wsLTW.Cells.Item(68, 12).Value = 2862.5
$wsLTW.Cells.Item(68, 14).Value = -4360.5

# LTW!row71
$wsLTW.Cells.Item(71, 8).Value = 2639.3462
$wsLTW.Cells.Item(71, 10).Value = 2862.5
$wsLTW.Cells.Item(71, 12).Value = 14312.5
$wsLTW.Cells.Item(71, 14).Value = -21800.5

# WVR!row122
$wsWVR.Cells.Item(122, 8).Value = 52632572
$wsWVR.Cells.Item(122, 9).Value = 71429540
$wsWVR.Cells.Item(122, 10).Value = 1076.6
$wsWVR.Cells.Item(122, 11).Value = 214288620
$wsWVR.Cells.Item(122, 12).Value = 3229.8
$wsWVR.Cells.Item(122, 13).Value = -214286170
$wsWVR.Cells.Item(122, 14).Value = -8129.799999999999

# WVR!row128
$wsWVR.Cells.Item(128, 8).Value = 52450
$wsWVR.Cells.Item(128, 10).Value = 52450
$wsWVR.Cells.Item(128, 12).Value = 52450
$wsWVR.Cells.Item(128, 14).Value = -62410
